$d = $word.ActiveDocument

# 1) Merge the split title runs "Engenharia de Software -  Cap. " + "2"
#    into a single run "Engenharia de Software -  Cap. 2".
$d.Content.Find.Execute("Engenharia de Software -  Cap. 2", $true, $false, $false, $false, $false, $true, 1, $false, "Engenharia de Software -  Cap. 2", 2)

# 2) Remove the word "não " from the Cascata description and split that
#    sentence's run in two, with the (moved) _GoBack bookmark sitting
#    between "anterior " and "foi concluído." - mirroring a manual edit
#    at that cursor position.
$d.Content.Find.Execute("não foi concluído", $true, $false, $false, $false, $false, $true, 1, $false, "foi concluído", 2)

# 3) Move the "_GoBack" bookmark (Word's "last edit location" marker) from
#    the end of the document to the point right after "anterior ".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("anterior ")
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)
